$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Item(1).TextToDisplay = "04a7ae36-c608-4111-b7b5-711d2e1f6785.md"
$ws1.Range("D2").Value2 = "2016-03-21 08:53:04"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Item(1).TextToDisplay = "04a7ae36-c608-4111-b7b5-711d2e1f6785.md"
$ws2.Hyperlinks.Item(2).TextToDisplay = "04a7ae36-c608-4111-b7b5-711d2e1f6785.df90fed8a2775395e0722b16788ca40b49d7b5f0.zh-cn.xlf"
$ws2.Range("E2").Value2 = "2016-03-21 08:53:00"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Item(1).TextToDisplay = "04a7ae36-c608-4111-b7b5-711d2e1f6785.md"
$ws3.Hyperlinks.Item(2).TextToDisplay = "04a7ae36-c608-4111-b7b5-711d2e1f6785.df90fed8a2775395e0722b16788ca40b49d7b5f0.de-de.xlf"
$ws3.Range("E2").Value2 = "2016-03-21 08:53:04"
